$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (230-233) appended after the existing data which ended at row 229
$data = @(
    @(44304, 3, 15, 182.3043266893534),
    @(44305, 1, 13, 157.997083130773),
    @(44306, 2, 12, 145.8434613514827),
    @(44307, 1, 13, 157.997083130773)
)

$startRow = 230
$endRow = $startRow + $data.Length - 1

# Copy the formatting (style) of the last existing data row (229) down onto the
# new rows A230:A233, matching the s="2" style used for the date column.
$srcA = $ws.Cells.Item($startRow - 1, 1)
$dstA = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$srcA.Copy() | Out-Null
$dstA.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
